$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.482.08'
$ws.Range('E2').Value = '  -0.89%  '
$ws.Range('D3').Value = '3.903.95'
$ws.Range('E3').Value = '  +3.91%  '
$ws.Range('E4').Value = '  -0.03%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.40'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +0.17%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.87'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -0.91%  '
$ws.Range('D7').Value = '3.906.49'
$ws.Range('E7').Value = '  +4.06%  '
$ws.Range('E8').Value = '  +0.02%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.526'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  -2.27%  '
$ws.Range('E10').Value = '  -3.84%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.37'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  -0.01%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.457'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  -0.35%  '
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.74'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  -2.53%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000244'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  -1.10%  '
$ws.Range('D15').Value = '4.565.43'
$ws.Range('E15').Value = '  +4.17%  '
$ws.Range('D16').Value = '3.900.59'
$ws.Range('E16').Value = '  +4.00%  '
$ws.Range('D17').Value = '68.717.25'
$ws.Range('E17').Value = '  -0.51%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.39'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('E19').Value = '  -1.04%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.97'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  -4.04%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.17'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  -0.54%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '482.82'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  -1.67%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.715'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  -1.59%  '
$ws.Range('E24').Value = '  +13.55%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.30'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  -0.42%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.22'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  -1.93%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.96'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  -2.55%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('E29').Value = '  -0.07%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.93'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  -1.05%  '
$ws.Range('D31').Value = '4.061.58'
$ws.Range('E31').Value = '  +4.08%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.83'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  -3.24%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.36'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  -2.80%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.87'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +0.48%  '
$ws.Range('D35').Value = '3.854.63'
$ws.Range('E35').Value = '  +3.61%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.106'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  -1.21%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.03'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  +2.68%  '
$ws.Range('E38').Value = '  -0.07%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.85'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  -1.24%  '
$ws.Range('E40').Value = '  +0.08%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.02'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  -2.33%  '
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.316'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  -2.36%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '432.52'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  +1.00%  '
$ws.Range('E44').Value = '  -0.17%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.97'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  -1.03%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.41'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -0.40%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '141.87'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  -0.55%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.816.65'
$ws.Range('E49').Value = '  +0.16%  '
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.93'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  +9.34%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '39.14'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  -2.76%  '
